$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.819.31'
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Value = '2.924.07'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '501.56'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").Value = '132.58'
$ws.Range("E6").Value = '  -4.67%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -3.19%  '
$ws.Range("D9").Value = '7.13'
$ws.Range("E9").Value = '  -5.22%  '
$ws.Range("E10").Value = '  -5.17%  '
$ws.Range("E11").Value = '  -4.07%  '
$ws.Range("D12").Value = '3.418.78'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("E13").Value = '  -4.05%  '
$ws.Range("D14").Value = '25.66'
$ws.Range("E14").Value = '  -2.91%  '
$ws.Range("D15").Value = '0.0000160'
$ws.Range("E15").Value = '  -2.65%  '
$ws.Range("D16").Value = '55.757.21'
$ws.Range("E16").Value = '  -2.83%  '
$ws.Range("D17").Value = '5.94'
$ws.Range("E17").Value = '  -4.47%  '
$ws.Range("D18").Value = '2.920.72'
$ws.Range("E18").Value = '  -2.99%  '
$ws.Range("D19").Value = '12.60'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("D21").Value = '313.97'
$ws.Range("E21").Value = '  -4.56%  '
$ws.Range("E23").Value = '  -2.45%  '
$ws.Range("D24").Value = '62.86'
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").Value = '3.036.97'
$ws.Range("E25").Value = '  -2.96%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -4.91%  '
$ws.Range("D28").Value = '0.0₃0837'
$ws.Range("E28").Value = '  -8.53%  '
$ws.Range("E29").Value = '  -7.02%  '
$ws.Range("D30").Value = '6.84'
$ws.Range("E30").Value = '  -8.27%  '
$ws.Range("E31").Value = '  -3.14%  '
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("D33").Value = '19.89'
$ws.Range("D34").Value = '150.99'
$ws.Range("E34").Value = '  -2.10%  '
$ws.Range("E35").Value = '  -7.31%  '
$ws.Range("D36").Value = '5.59'
$ws.Range("E36").Value = '  -4.84%  '
$ws.Range("D37").Value = '23.94'
$ws.Range("E37").Value = '  -2.15%  '
$ws.Range("E38").Value = '  -7.25%  '
$ws.Range("D39").Value = '0.0645'
$ws.Range("E39").Value = '  -5.22%  '
$ws.Range("D40").Value = '36.40'
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -3.48%  '
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").Value = '2.114.12'
$ws.Range("E44").Value = '  -7.89%  '
$ws.Range("D45").Value = '5.96'
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("E46").Value = '  -5.65%  '
$ws.Range("D47").Value = '0.919'
$ws.Range("E47").Value = '  -6.78%  '
$ws.Range("E48").Value = '  -2.74%  '
$ws.Range("D49").Value = '18.65'
$ws.Range("E49").Value = '  -4.07%  '
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  -9.32%  '
